# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# --- Rushing sheet updates ---
$wsRushing = $wb.Worksheets.Item("Rushing")

$wsRushing.Range("C2").Value = 20
$wsRushing.Range("D2").Value = 23
$wsRushing.Range("E2").Value = 20

$wsRushing.Range("C5").Value = 117
$wsRushing.Range("D5").Value = 75
$wsRushing.Range("F5").Value = 34

$wsRushing.Range("C6").Value = 24

$wsRushing.Range("D11").Value = 1

$wsRushing.Range("D13").Value = 8

# --- Receiving sheet updates ---
$wsReceiving = $wb.Worksheets.Item("Receiving")

$wsReceiving.Range("C2").Value = 58
$wsReceiving.Range("D2").Value = 41
$wsReceiving.Range("E2").Value = 7
$wsReceiving.Range("F2").Value = 4
$wsReceiving.Range("G2").Value = 14
$wsReceiving.Range("H2").Value = 9

$wsReceiving.Range("C8").Value = 2
$wsReceiving.Range("D8").Value = 2

$wsReceiving.Range("C10").Value = 58
$wsReceiving.Range("D10").Value = 36
$wsReceiving.Range("E10").Value = 25
$wsReceiving.Range("F10").Value = 9
$wsReceiving.Range("G10").Value = 8

$wsReceiving.Range("C11").Value = 42
$wsReceiving.Range("D11").Value = 28

$wsReceiving.Range("C13").Value = 10
$wsReceiving.Range("D13").Value = 6

$wsReceiving.Range("C20").Value = 14
$wsReceiving.Range("D20").Value = 9
$wsReceiving.Range("G20").Value = 6
$wsReceiving.Range("H20").Value = 4
